$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.007.80'
$ws.Range('E2').Value = '  -1.74%  '
$ws.Range('D3').Value = '1.908.55'
$ws.Range('E3').Value = '  -3.22%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = '324.80'
$ws.Range('E5').Value = '  -0.75%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('E7').Value = '  -1.56%  '
$ws.Range('D8').Value = '0.3824'
$ws.Range('E8').Value = '  -2.27%  '
$ws.Range('D9').Value = '0.07715'
$ws.Range('E9').Value = '  -2.79%  '
$ws.Range('D10').Value = '0.9800'
$ws.Range('E10').Value = '  -0.82%  '
$ws.Range('D11').Value = '22.09'
$ws.Range('E11').Value = '  -2.90%  '
$ws.Range('D12').Value = '1.892.46'
$ws.Range('E12').Value = '  -6.11%  '
$ws.Range('D13').Value = '5.674'
$ws.Range('E13').Value = '  -2.21%  '
$ws.Range('D14').Value = '6.934'
$ws.Range('E14').Value = '  -3.41%  '
$ws.Range('D15').Value = '0.07040'
$ws.Range('E15').Value = '  -1.22%  '
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').Value = '83.82'
$ws.Range('E17').Value = '  -4.50%  '
$ws.Range('D18').Value = '0.000009469'
$ws.Range('E18').Value = '  -4.53%  '
$ws.Range('D19').Value = '16.68'
$ws.Range('E19').Value = '  -3.21%  '
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.57%  '
$ws.Range('D21').Value = '28.983.83'
$ws.Range('E21').Value = '  -2.05%  '
$ws.Range('D22').Value = '5.323'
$ws.Range('E22').Value = '  -4.00%  '
$ws.Range('E23').Value = '  -2.26%  '
$ws.Range('D24').Value = '2.126.93'
$ws.Range('E24').Value = '  -5.42%  '
$ws.Range('D25').Value = '2.093'
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('D26').Value = '158.27'
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('D27').Value = '19.06'
$ws.Range('E27').Value = '  -2.45%  '
$ws.Range('D28').Value = '5.659'
$ws.Range('E28').Value = '  -2.33%  '
$ws.Range('D29').Value = '117.48'
$ws.Range('E29').Value = '  -1.88%  '
$ws.Range('D30').Value = '1.855'
$ws.Range('E30').Value = '  -1.78%  '
$ws.Range('D31').Value = '0.09285'
$ws.Range('E31').Value = '  -1.49%  '
$ws.Range('D32').Value = '0.8656'
$ws.Range('E32').Value = '  -1.45%  '
$ws.Range('D33').Value = '5.069'
$ws.Range('E33').Value = '  -3.14%  '
$ws.Range('D34').Value = '1.246'
$ws.Range('E34').Value = '  -5.30%  '
$ws.Range('D35').Value = '3.019'
$ws.Range('E35').Value = '  -4.20%  '
$ws.Range('D36').Value = '0.05733'
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('E37').Value = '  -0.98%  '
$ws.Range('D38').Value = '1.001'
$ws.Range('E38').Value = '  -0.40%  '
$ws.Range('D39').Value = '0.02037'
$ws.Range('E39').Value = '  -3.16%  '
$ws.Range('D40').Value = '0.5509'
$ws.Range('E40').Value = '  -3.44%  '
$ws.Range('D41').Value = '7.405'
$ws.Range('E41').Value = '  -3.88%  '
$ws.Range('E42').Value = '  -2.24%  '
$ws.Range('D43').Value = '2.862'
$ws.Range('E43').Value = '  +3.77%  '
$ws.Range('D44').Value = '9.350'
$ws.Range('E44').Value = '  -2.86%  '
$ws.Range('D45').Value = '0.5187'
$ws.Range('E45').Value = '  -2.67%  '
$ws.Range('D46').Value = '11.21'
$ws.Range('E46').Value = '  -4.79%  '
$ws.Range('D47').Value = '0.06835'
$ws.Range('E47').Value = '  -1.40%  '
$ws.Range('D48').Value = '2.048'
$ws.Range('E48').Value = '  -4.37%  '
$ws.Range('D49').Value = '110.99'
$ws.Range('E49').Value = '  -2.15%  '
$ws.Range('D50').Value = '1.780'
$ws.Range('E50').Value = '  -2.50%  '
$ws.Range('D51').Value = '0.000002556'
$ws.Range('E51').Value = '  -6.51%  '
